# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (standard Office palette) - only
#                           wired to the (unused) notes master.
#   ppt/theme/theme2.xml -> "Integral" (green/yellow palette) - the theme
#                           actually used by the presentation's slide master
#                           / all slides.
# The authored edit swaps the content of the two theme parts: the deck's
# active theme becomes the plain "Office Theme" palette (and the Integral
# palette ends up parked in the no-longer-visible notes-master theme slot).
#
# The PowerPoint object model's ColorScheme exposes exactly the 12 colour
# slots of the active theme's <a:clrScheme> (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink), in that order, and writing RGB there edits the theme
# part backing the slide master (the presentation's real, visible theme).
# Push the standard Office palette into every slot to perform the swap.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

$officeRGB = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $officeRGB.Count; $i++) {
    $cs.Item($i).RGB = $officeRGB[$i - 1]
}
